$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update Version value (B3): 0.1.1 -> 0.2.0
$ws.Cells.Item(3, 2).Value = "0.2.0"

# 2. Update Date value (B8): 2023-10-19T16:17:18+00:00 -> 2023-10-19T17:05:12+00:00
$ws.Cells.Item(8, 2).Value = "2023-10-19T17:05:12+00:00"

# 3. Insert a new "Jurisdiction" / "iso:code:3166:FR" row at row 11, pushing the
#    existing Description/Purpose/Copyright/Immutable rows down by one
#    (rows 11-14 -> 12-15).
#
#    Capture the current values first (top to bottom) so we know exactly what
#    to re-write, including blanks, into the shifted rows.
$srcA = @{}
$srcB = @{}
for ($r = 11; $r -le 14; $r++) {
    $srcA[$r] = $ws.Cells.Item($r, 1).Value2
    $srcB[$r] = $ws.Cells.Item($r, 2).Value2
}

# Shift bottom-up using cell-to-cell Copy so the existing cell style (s=2) is
# reused instead of a brand new style entry being generated, then explicitly
# (re)write the value - including clearing to blank where the source was
# blank, since Copy() from a blank cell does not clear a non-blank target.
for ($r = 14; $r -ge 11; $r--) {
    $ws.Cells.Item($r, 1).Copy($ws.Cells.Item($r + 1, 1))
    $ws.Cells.Item($r + 1, 1).Value = $srcA[$r]

    $ws.Cells.Item($r, 2).Copy($ws.Cells.Item($r + 1, 2))
    $ws.Cells.Item($r + 1, 2).Value = $srcB[$r]
}

# Now write the new Jurisdiction row into row 11 (style already copied above)
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "iso:code:3166:FR"
